$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relocate the "primary key / foreign key" legend entries in columns J:K ---
# Move J4:K4 ("기본키") up to J3:K3
$ws.Range("J4:K4").Cut($ws.Range("J3"))
# Move J6:K6 ("외래키") up to J5:K5
$ws.Range("J6:K6").Cut($ws.Range("J5"))

# Clear the now-vacated cells and restore their plain formatting
$ws.Range("J4:K4").Clear()
$ws.Range("J6:K6").Clear()
$ws.Range("I4").Copy()
$ws.Range("J4:K4").PasteSpecial(-4122)

# --- Add a new "범례" (Legend) title above the legend box ---
$ws.Range("J2:K2").Merge()
$ws.Range("J2").Value = "범례"
$ws.Range("J2:K2").HorizontalAlignment = -4108

# --- Remove two blank rows from the lower empty block, shifting the ---
# --- trailing (내용/컬럼명/타입) labels up from rows 18-20 to rows 16-18 ---
$ws.Range("A16:K17").Delete(-4162)

# Restore the selection to match the edited workbook's last position
$ws.Range("F18").Select()
